# Student entry using excel
# The module to receive student data using excel done
#
# - Repurpose the "Guardian Contact" column (I) into a "Class Name" column
# - Update the sample row accordingly (Guardian phone number -> sample class name)
# - Add a second sample student row (row 3) showing the new Class Name field
# - Clear the stray custom number-format that had been applied to the
#   "Year" column (E) so it reverts to General

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-purpose column I: "Guardian Contact" -> "Class Name" ---
$ws.Range("I1").Value = "Class Name"
$ws.Range("I2").Value = "Science 1"

# --- Second sample student row ---
$ws.Range("B3").Value = "Amega"
$ws.Range("C3").Value = "Aisha"
$ws.Range("D3").Value = "Female"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "House three"
$ws.Range("G3").Value = "Day"
$ws.Range("H3").Value = "Visual Arts"
$ws.Range("I3").Value = "Vis 2"

# New row's text cells match the existing table's text-style formatting
$ws.Range("B3:D3").NumberFormat = "@"
$ws.Range("F3:I3").NumberFormat = "@"

# --- Drop the one-off number format that column E (Year) had ---
$ws.Columns.Item(5).ClearFormats()

# --- Column I should size itself to the new "Class Name" content ---
$ws.Columns.Item(9).AutoFit()

# --- Selection moves onto the (now full-height) Class Name column ---
$ws.Range("I1:I1048576").Select()
